$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header row (row 1) with new short column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case municipality / state names (apply .title()-style capitalization
#    to Spanish connector words like "de", "del", "el", "y", etc.)
$ws.Range('B5').Value = 'Pabellón De Arteaga'
$ws.Range('B6').Value = 'Rincón De Romos'
$ws.Range('B10').Value = 'Playas De Rosarito'
$ws.Range('B27').Value = 'Amatenango De La Frontera'
$ws.Range('B30').Value = 'Bejucal De Ocampo'
$ws.Range('B32').Value = 'Benemérito De Las Américas'
$ws.Range('B37').Value = 'Chiapa De Corzo'
$ws.Range('B41').Value = 'Comitán De Domínguez'
$ws.Range('B58').Value = 'Mazapa De Madero'
$ws.Range('B70').Value = 'Salto De Agua'
$ws.Range('B71').Value = 'San Cristóbal De Las Casas'
$ws.Range('B99').Value = 'Hidalgo Del Parral'
$ws.Range('B115').Value = 'San Juan De Sabinas'
$ws.Range('A123').Value = 'Ciudad De México'
$ws.Range('B127').Value = 'Cuajimalpa De Morelos'
$ws.Range('B149').Value = 'San Juan De Guadalupe'
$ws.Range('B150').Value = 'San Juan Del Río'
$ws.Range('A154').Value = 'Estado De México'
$ws.Range('B154').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B157').Value = 'Almoloya De Alquisiras'
$ws.Range('B158').Value = 'Almoloya De Juárez'
$ws.Range('B161').Value = 'Atizapán De Zaragoza'
$ws.Range('B174').Value = 'Ecatepec De Morelos'
$ws.Range('B177').Value = 'Ixtapan De La Sal'
$ws.Range('B185').Value = 'Naucalpan De Juárez'
$ws.Range('B194').Value = 'San Felipe Del Progreso'
$ws.Range('B195').Value = 'San Simón De Guerrero'
$ws.Range('B196').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B204').Value = 'Tenango Del Valle'
$ws.Range('B212').Value = 'Tlalnepantla De Baz'
$ws.Range('B216').Value = 'Villa De Allende'
$ws.Range('B217').Value = 'Villa Del Carbón'
$ws.Range('B227').Value = 'San Miguel De Allende'
$ws.Range('B228').Value = 'Apaseo El Alto'
$ws.Range('B229').Value = 'Apaseo El Grande'
$ws.Range('B236').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B240').Value = 'Jaral Del Progreso'
$ws.Range('B250').Value = 'San Diego De La Unión'
$ws.Range('B252').Value = 'San Francisco Del Rincón'
$ws.Range('B254').Value = 'San Luis De La Paz'
$ws.Range('B256').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B257').Value = 'Silao De La Victoria'
$ws.Range('B261').Value = 'Valle De Santiago'
$ws.Range('B266').Value = 'Acapulco De Juárez'
$ws.Range('B268').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B269').Value = 'Alcozauca De Guerrero'
$ws.Range('B273').Value = 'Atenango Del Río'
$ws.Range('B275').Value = 'Atoyac De Álvarez'
$ws.Range('B276').Value = 'Ayutla De Los Libres'
$ws.Range('B279').Value = 'Chilapa De Álvarez'
$ws.Range('B280').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B281').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B284').Value = 'Coyuca De Benítez'
$ws.Range('B285').Value = 'Coyuca De Catalán'
$ws.Range('B289').Value = 'Cutzamala De Pinzón'
$ws.Range('B295').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B296').Value = 'Iguala De La Independencia'
$ws.Range('B298').Value = 'Zihuatanejo De Azueta'
$ws.Range('B300').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B314').Value = 'Taxco De Alarcón'
$ws.Range('B316').Value = 'Técpan De Galeana'
$ws.Range('B318').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B320').Value = 'Tixtla De Guerrero'
$ws.Range('B323').Value = 'Tlapa De Comonfort'
$ws.Range('B335').Value = 'Agua Blanca De Iturbide'
$ws.Range('B340').Value = 'Atotonilco El Grande'
$ws.Range('B345').Value = 'Cuautepec De Hinojosa'
$ws.Range('B349').Value = 'Huasca De Ocampo'
$ws.Range('B352').Value = 'Huejutla De Reyes'
$ws.Range('B355').Value = 'Jacala De Ledezma'
$ws.Range('B361').Value = 'Mineral Del Chico'
$ws.Range('B362').Value = 'Mineral Del Monte'
$ws.Range('B363').Value = 'Mixquiahuala De Juárez'
$ws.Range('B364').Value = 'Molango De Escamilla'
$ws.Range('B366').Value = 'Omitlán De Juárez'
$ws.Range('B367').Value = 'Pachuca De Soto'
$ws.Range('B370').Value = 'Progreso De Obregón'
$ws.Range('B375').Value = 'Santiago De Anaya'
$ws.Range('B376').Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range('B380').Value = 'Tenango De Doria'
$ws.Range('B382').Value = 'Tepehuacán De Guerrero'
$ws.Range('B383').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B384').Value = 'Tezontepec De Aldama'
$ws.Range('B390').Value = 'Tula De Allende'
$ws.Range('B391').Value = 'Tulancingo De Bravo'
$ws.Range('B394').Value = 'Zacualtipán De Ángeles'
$ws.Range('B401').Value = 'Autlán De Navarro'
$ws.Range('B411').Value = 'Encarnación De Díaz'
$ws.Range('B415').Value = 'Ixtlahuacán Del Río'
$ws.Range('B419').Value = 'Lagos De Moreno'
$ws.Range('B425').Value = 'San Cristóbal De La Barranca'
$ws.Range('B426').Value = 'San Diego De Alejandría'
$ws.Range('B427').Value = 'San Juan De Los Lagos'
$ws.Range('B429').Value = 'San Miguel El Alto'
$ws.Range('B432').Value = 'Tamazula De Gordiano'
$ws.Range('B435').Value = 'Teocuitatlán De Corona'
$ws.Range('B436').Value = 'Tepatitlán De Morelos'
$ws.Range('B438').Value = 'Tizapán El Alto'
$ws.Range('B439').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B443').Value = 'Unión De Tula'
$ws.Range('B446').Value = 'Zapotlán El Grande'
$ws.Range('B507').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B528').Value = 'Coatlán Del Río'
$ws.Range('B538').Value = 'Puente De Ixtla'
$ws.Range('B542').Value = 'Tetela Del Volcán'
$ws.Range('B552').Value = 'Ixtlán Del Río'
$ws.Range('B568').Value = 'Mier Y Noriega'
$ws.Range('B571').Value = 'San Nicolás De Los Garza'
$ws.Range('B575').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B580').Value = 'Coicoyán De Las Flores'
$ws.Range('B582').Value = 'Constancia Del Rosario'
$ws.Range('B584').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B585').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B586').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B588').Value = 'Ixtlán De Juárez'
$ws.Range('B589').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B594').Value = 'Mariscala De Juárez'
$ws.Range('B595').Value = 'Mártires De Tacubaya'
$ws.Range('B598').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B601').Value = 'Nejapa De Madero'
$ws.Range('B602').Value = 'Oaxaca De Juárez'
$ws.Range('B603').Value = 'Ocotlán De Morelos'
$ws.Range('B604').Value = 'Pinotepa De Don Luis'
$ws.Range('B606').Value = 'Putla Villa De Guerrero'
$ws.Range('B610').Value = 'San Agustín De Las Juntas'
$ws.Range('B618').Value = 'San Dionisio Del Mar'
$ws.Range('B622').Value = 'San Felipe Jalapa De Díaz'
$ws.Range('B626').Value = 'San Francisco Del Mar'
$ws.Range('B632').Value = 'San José Del Progreso'
$ws.Range('B634').Value = 'San Juan Bautista Lo De Soto'
$ws.Range('B671').Value = 'San Miguel Del Puerto'
$ws.Range('B679').Value = 'San Pedro El Alto'
$ws.Range('B697').Value = 'Santa Ana Del Valle'
$ws.Range('B712').Value = 'Santa María Del Tule'
$ws.Range('B716').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B755').Value = 'Santo Domingo De Morelos'
$ws.Range('B767').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B768').Value = 'Tataltepec De Valdés'
$ws.Range('B769').Value = 'Teococuilco De Marcos Pérez'
$ws.Range('B770').Value = 'Teotitlán Del Valle'
$ws.Range('B772').Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range('B773').Value = 'Tlacolula De Matamoros'
$ws.Range('B774').Value = 'Totontepec Villa De Morelos'
$ws.Range('B777').Value = 'Villa De Tututepec'
$ws.Range('B778').Value = 'Villa De Zaachila'
$ws.Range('B780').Value = 'Villa Sola De Vega'
$ws.Range('B781').Value = 'Zimatlán De Álvarez'
$ws.Range('B792').Value = 'Ayotoxco De Guerrero'
$ws.Range('B795').Value = 'Chalchicomula De Sesma'
$ws.Range('B807').Value = 'Cuetzalan Del Progreso'
$ws.Range('B821').Value = 'Izúcar De Matamoros'
$ws.Range('B827').Value = 'Mazapiltepec De Juárez'
$ws.Range('B832').Value = 'Palmar De Bravo'
$ws.Range('B843').Value = 'San Nicolás De Los Ranchos'
$ws.Range('B844').Value = 'San Salvador El Seco'
$ws.Range('B845').Value = 'San Salvador El Verde'
$ws.Range('B847').Value = 'Tecali De Herrera'
$ws.Range('B853').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B858').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B871').Value = 'Xochitlán De Vicente Suárez'
$ws.Range('B881').Value = 'Amealco De Bonfil'
$ws.Range('B883').Value = 'Cadereyta De Montes'
$ws.Range('B887').Value = 'Jalpan De Serra'
$ws.Range('B888').Value = 'Landa De Matamoros'
$ws.Range('B890').Value = 'Pinal De Amoles'
$ws.Range('B892').Value = 'San Juan Del Río'
$ws.Range('B903').Value = 'Armadillo De Los Infante'
$ws.Range('B904').Value = 'Axtla De Terrazas'
$ws.Range('B910').Value = 'Ciudad Del Maíz'
$ws.Range('B919').Value = 'Mexquitic De Carmona'
$ws.Range('B922').Value = 'San Ciro De Acosta'
$ws.Range('B926').Value = 'Santa María Del Río'
$ws.Range('B928').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B936').Value = 'Villa De Guadalupe'
$ws.Range('B937').Value = 'Villa De La Paz'
$ws.Range('B938').Value = 'Villa De Ramos'
$ws.Range('B939').Value = 'Villa De Reyes'
$ws.Range('B970').Value = 'Jalpa De Méndez'
$ws.Range('B996').Value = 'Soto La Marina'
$ws.Range('B1007').Value = 'Contla De Juan Cuamatzi'
$ws.Range('B1010').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B1011').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B1013').Value = 'Papalotla De Xicohténcatl'
$ws.Range('B1014').Value = 'Sanctórum De Lázaro Cárdenas'
$ws.Range('B1026').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B1030').Value = 'Amatlán De Los Reyes'
$ws.Range('B1039').Value = 'Boca Del Río'
$ws.Range('B1041').Value = 'Camarón De Tejeda'
$ws.Range('B1045').Value = 'Cazones De Herrera'
$ws.Range('B1059').Value = 'Cosamaloapan De Carpio'
$ws.Range('B1060').Value = 'Cosautlán De Carvajal'
$ws.Range('B1074').Value = 'Hueyapan De Ocampo'
$ws.Range('B1075').Value = 'Ignacio De La Llave'
$ws.Range('B1078').Value = 'Ixhuacán De Los Reyes'
$ws.Range('B1079').Value = 'Ixhuatlán De Madero'
$ws.Range('B1080').Value = 'Ixhuatlán Del Café'
$ws.Range('B1081').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B1090').Value = 'Juchique De Ferrer'
$ws.Range('B1094').Value = 'Lerdo De Tejada'
$ws.Range('B1098').Value = 'Martínez De La Torre'
$ws.Range('B1101').Value = 'Medellín De Bravo'
$ws.Range('B1112').Value = 'Ozuluama De Mascareñas'
$ws.Range('B1116').Value = 'Paso De Ovejas'
$ws.Range('B1117').Value = 'Paso Del Macho'
$ws.Range('B1120').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1128').Value = 'Sayula De Alemán'
$ws.Range('B1132').Value = 'Soledad De Doblado'
$ws.Range('B1163').Value = 'Vega De Alatorre'
$ws.Range('B1173').Value = 'Zozocolco De Hidalgo'
$ws.Range('B1192').Value = 'Moyahua De Estrada'
$ws.Range('B1193').Value = 'Noria De Ángeles'
$ws.Range('B1201').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1203').Value = 'Villa De Cos'

# 3. Correct tiny floating-point rounding differences in column D
$ws.Range('D9').Value = 0.0009941593140300731
$ws.Range('D30').Value = 0.0009941593140300731
$ws.Range('D31').Value = 0.0009941593140300731
$ws.Range('D33').Value = 0.0009941593140300731
$ws.Range('D82').Value = 0.0009941593140300731
$ws.Range('D191').Value = 0.0009941593140300731
$ws.Range('D218').Value = 0.0009941593140300731
$ws.Range('D371').Value = 0.0009941593140300731
$ws.Range('D380').Value = 0.0009941593140300731
$ws.Range('D562').Value = 0.0009941593140300731
$ws.Range('D588').Value = 0.0009941593140300731
$ws.Range('D598').Value = 0.0009941593140300731
$ws.Range('D701').Value = 0.0009941593140300731
$ws.Range('D744').Value = 0.0009941593140300731
$ws.Range('D821').Value = 0.0009941593140300731
$ws.Range('D874').Value = 0.0009941593140300731
$ws.Range('D901').Value = 0.0009941593140300731
$ws.Range('D919').Value = 0.0009941593140300731
$ws.Range('D1080').Value = 0.0009941593140300731
$ws.Range('D1094').Value = 0.0009941593140300731
$ws.Range('D1114').Value = 0.0009941593140300731

# 4. Remove trailing footer/metadata rows (1211-1216), shrinking the sheet to A1:D1210
$ws.Range("A1211:D1216").EntireRow.Delete()

